# edit.ps1 - apply "Bad Ideas" slide content revisions
# Target slide is the only slide in the deck (the "What not to do" /
# "Week before final" / "During Final" triptych).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape 3: "What not to do" -> "Bad Ideas (these have all happened)" box
# ---------------------------------------------------------------------
$shpBad = $s.Shapes.Item(3)
$trBad = $shpBad.TextFrame.TextRange

# Paragraph 1 title: split into three runs ("Bad " / "Ideas (these " / "have all happened)")
$title = $trBad.Paragraphs(1, 1)
$title.Runs(1, 1).Text = "Bad Ideas (these have all happened)"
$titleMid = $title.Characters(5, 13)          # "Ideas (these "
$titleMid.Font.Underline = 1                   # force run split, keep underline formatting

# Paragraph 2 bullet
$trBad.Paragraphs(2, 1).Runs(1, 1).Text = "Skip the final"

# Paragraph 4 bullet
$trBad.Paragraphs(4, 1).Runs(1, 1).Text = "Show code instead of diagrams"

# Paragraph 5 bullet
$trBad.Paragraphs(5, 1).Runs(1, 1).Text = "Video > 2m or < 1:30m"

# Paragraph 7 bullet
$trBad.Paragraphs(7, 1).Runs(1, 1).Text = "Play from YouTube with sketchy suggested videos"

# Paragraph 8 bullet
$trBad.Paragraphs(8, 1).Runs(1, 1).Text = "Skip practicing"

# Paragraph 9 bullet
$trBad.Paragraphs(9, 1).Runs(1, 1).Text = "Bring up video with lengthy login, download, laptop, AV connection, USB connection, etc."

# Paragraphs 12 & 13 swap content
$trBad.Paragraphs(12, 1).Runs(1, 1).Text = "Complain about how bad your project is "
$trBad.Paragraphs(13, 1).Runs(1, 1).Text = "Talk to the projector screen"

# ---------------------------------------------------------------------
# Shape 4: "Week before final:" box
# ---------------------------------------------------------------------
$shpWeek = $s.Shapes.Item(4)
$trWeek = $shpWeek.TextFrame.TextRange

# Paragraph 3 bullet
$trWeek.Paragraphs(3, 1).Runs(1, 1).Text = "Combine into one video, 1:30m to 2:00m, w/o audio"

# Paragraph 4 bullet: "Upload to Canvas, "UNO Academic Video" (VidGrid)"
#                  -> "Upload to Canvas (VigGrid or VuJa)"
$uploadPara = $trWeek.Paragraphs(4, 1)
$uploadPara.Runs(1, 1).Text = "Upload to Canvas ("
$uploadPara.Runs(2, 1).Text = "VigGrid or VuJa"

$uploadText = $uploadPara.Text
$base = $uploadText.IndexOf("VigGrid or VuJa") + 1

$vigRange = $uploadPara.Characters($base, 7)
$vigRange.Text = "VigGrid"

$orRange = $uploadPara.Characters($base + 7, 4)
$orRange.Text = " or "

$vujaRange = $uploadPara.Characters($base + 11, 4)
$vujaRange.Text = "VuJa"

# ---------------------------------------------------------------------
# Shape 4 Paragraph 5 bullet: italic "different computer"
# ---------------------------------------------------------------------
$copyPara = $trWeek.Paragraphs(5, 1)
$copyPara.Runs(1, 1).Text = "Copy minified link and check on a different computer"
$prefix = "Copy minified link and check on a "
$italicRange = $copyPara.Characters($prefix.Length + 1, ("different computer").Length)
$italicRange.Font.Italic = -1

# ---------------------------------------------------------------------
# Shape 5: "During Final:" box
# ---------------------------------------------------------------------
$shpDuring = $s.Shapes.Item(5)
$trDuring = $shpDuring.TextFrame.TextRange

$trDuring.Paragraphs(2, 1).Runs(1, 1).Text = "Bring up your video with written, minified link"
$trDuring.Paragraphs(3, 1).Runs(1, 1).Text = "Discuss main points and everything else cool"
$trDuring.Paragraphs(4, 1).Runs(1, 1).Text = "Be positive and honest: Talk to future employer"
